# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.652.94'
$ws.Range("D3").Value = '2.386.13'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '504.23'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.31'
$ws.Range("E6").Value = '  +2.23%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.549'
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("D9").Value = '2.388.43'
$ws.Range("E9").Value = '  -0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0972'
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.151'
$ws.Range("E11").Value = '  +0.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.324'
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.64'
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").Value = '2.808.17'
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").Value = '56.571.34'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.62'
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '2.387.90'
$ws.Range("E18").Value = '  +2.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.17'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '308.95'
$ws.Range("E21").Value = '  -1.05%  '
$ws.Range("E22").Value = '  +1.14%  '
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("E24").Value = '  -4.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.94'
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.380'
$ws.Range("E27").Value = '  +2.78%  '
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.40'
$ws.Range("E29").Value = '  +2.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '176.39'
$ws.Range("E30").Value = '  +1.25%  '
$ws.Range("D31").Value = '0.0₃0725'
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("E32").Value = '  -0.47%  '
$ws.Range("E33").Value = '  +1.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.84'
$ws.Range("E34").Value = '  -4.34%  '
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.78'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.20'
$ws.Range("E38").Value = '  -2.60%  '
$ws.Range("E39").Value = '  +1.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.73'
$ws.Range("E40").Value = '  +2.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.822'
$ws.Range("E41").Value = '  +5.81%  '
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '130.69'
$ws.Range("E43").Value = '  +0.58%  '
$ws.Range("E44").Value = '  +0.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.82'
$ws.Range("E45").Value = '  -1.13%  '
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("E47").Value = '  +1.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '247.62'
$ws.Range("E48").Value = '  -2.41%  '
$ws.Range("E49").Value = '  -0.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0209'
$ws.Range("E50").Value = '  +1.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.17'
$ws.Range("E51").Value = '  +6.61%  '
